$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D from automatic numeric/date conversion while we set text-like values
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '49.429.30'
$ws.Range("E2").Value = '  -1.07%  '

$ws.Range("D3").Value = '2.629.15'
$ws.Range("E3").Value = '  -1.00%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = '111.72'
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").Value = '324.52'
$ws.Range("E6").Value = '  -1.08%  '

$ws.Range("D7").Value = '0.525'
$ws.Range("E7").Value = '  -1.26%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").Value = '0.543'
$ws.Range("E9").Value = '  -3.12%  '

$ws.Range("D10").Value = '39.54'
$ws.Range("E10").Value = '  -3.24%  '

$ws.Range("E11").Value = '  -3.75%  '

$ws.Range("D12").Value = '0.0810'
$ws.Range("E12").Value = '  -1.71%  '

$ws.Range("E13").Value = '  +1.37%  '

$ws.Range("E14").Value = '  -0.58%  '

$ws.Range("D15").Value = '3.047.26'
$ws.Range("E15").Value = '  -0.78%  '

$ws.Range("D16").Value = '2.638.41'
$ws.Range("E16").Value = '  +1.76%  '

$ws.Range("D17").Value = '0.847'
$ws.Range("E17").Value = '  -4.04%  '

$ws.Range("D18").Value = '49.388.11'
$ws.Range("E18").Value = '  -1.10%  '

$ws.Range("D19").Value = '12.85'
$ws.Range("E19").Value = '  -3.55%  '

$ws.Range("E20").Value = '  -1.24%  '

$ws.Range("D21").Value = '6.68'
$ws.Range("E21").Value = '  -2.33%  '

$ws.Range("D22").Value = '0.0₃0945'
$ws.Range("E22").Value = '  -2.07%  '

$ws.Range("D23").Value = '269.69'
$ws.Range("E23").Value = '  -4.08%  '

$ws.Range("D24").Value = '68.83'
$ws.Range("E24").Value = '  -5.98%  '

$ws.Range("E25").Value = '  -2.78%  '

$ws.Range("D26").Value = '26.18'
$ws.Range("E26").Value = '  -3.34%  '

$ws.Range("E27").Value = '  +0.06%  '

$ws.Range("E28").Value = '  +2.74%  '

$ws.Range("E29").Value = '  -1.45%  '

$ws.Range("E30").Value = '  -4.97%  '

$ws.Range("D31").Value = '34.61'
$ws.Range("E31").Value = '  -6.62%  '

$ws.Range("D32").Value = '49.44'
$ws.Range("E32").Value = '  -0.69%  '

$ws.Range("D33").Value = '5.47'
$ws.Range("E33").Value = '  +0.36%  '

$ws.Range("D34").Value = '0.0812'
$ws.Range("E34").Value = '  +1.32%  '

$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("D36").Value = '18.83'
$ws.Range("E36").Value = '  -4.12%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = '2.04'
$ws.Range("E37").Value = '  -1.87%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '4.87'
$ws.Range("E38").Value = '  +1.74%  '

$ws.Range("D39").Value = '3.09'
$ws.Range("E39").Value = '  -1.32%  '

$ws.Range("D40").Value = '128.89'
$ws.Range("E40").Value = '  +1.06%  '

$ws.Range("E41").Value = '  -1.96%  '

$ws.Range("D42").Value = '22.04'
$ws.Range("E42").Value = '  -1.82%  '

$ws.Range("E43").Value = '  +3.92%  '

$ws.Range("E44").Value = '  -3.61%  '

$ws.Range("D45").Value = '2.055.78'
$ws.Range("E45").Value = '  -0.72%  '

$ws.Range("D46").Value = '3.20'
$ws.Range("E46").Value = '  -5.50%  '

$ws.Range("D47").Value = '2.11'
$ws.Range("E47").Value = '  +5.97%  '

$ws.Range("E48").Value = '  -5.43%  '

$ws.Range("D49").Value = '8.90'
$ws.Range("E49").Value = '  -2.14%  '

$ws.Range("D50").Value = '59.08'
$ws.Range("E50").Value = '  +1.80%  '

$ws.Range("E51").Value = '  -4.46%  '

# Restore default (General) style/number formatting on column D so the sheet matches the source styling
$ws.Range("D2:D51").Style = "Normal"

